$p = $ppt.ActivePresentation
Write-Output ("Before Designs.Count: " + $p.Designs.Count)
try {
  $nd = $p.Designs.Add()
  Write-Output ("Added: " + $nd)
  Write-Output ("After Designs.Count: " + $p.Designs.Count)
} catch {
  Write-Output ("ERR: " + $_.Exception.Message)
}
